# Update R training schedule: swap paired webinar topics/links (moves
# Tidyverse Introduction up to week 4), then extend the sheet with 10
# additional blank (date-formatted) rows for future weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Swap-WebinarRows($r1, $r2) {
    $e1 = $ws.Range("E$r1").Value()
    $h1 = $ws.Range("H$r1").Value()
    $e2 = $ws.Range("E$r2").Value()
    $h2 = $ws.Range("H$r2").Value()
    $ws.Range("E$r1").Value = $e2
    $ws.Range("H$r1").Value = $h2
    $ws.Range("E$r2").Value = $e1
    $ws.Range("H$r2").Value = $h1
}

Swap-WebinarRows 5 6
Swap-WebinarRows 7 8
Swap-WebinarRows 9 10
Swap-WebinarRows 11 12

# Extend the schedule with 10 more (currently blank) week rows, carrying
# over the date number-format used by the existing week_start/week_end
# columns (B:C).
$fmtSrc = $ws.Range("B2:C2")
$newRows = $ws.Range("B20:C29")
$fmtSrc.Copy($newRows)
$newRows.ClearContents()

$ws.Range("E13").Select()
